$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1190320826869504, 0.306821227259698, 3.537761648806719, 0.4942365360607697)
    3  = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 10.19245300693656)
    4  = @(0.2917716402565462, 0.306821227259698, 0.7527432677738641, 0.4942365360607697)
    5  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    6  = @(0.2917716402565462, 0.306821227259698, 0.7527432677738641, 0.4942365360607697)
    7  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697)
    8  = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697)
    9  = @(0.04271373187048222, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    10 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]
    $g = $b + $c + $d + $e

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = $g
}
